$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checks")
$ws.Range("D10").Value = "Rule4"
